$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Silver Rear_side")
$ws1.Range("B10").NumberFormat = "@"
$ws1.Range("B10").Value = "5,221"

$ws2 = $wb.Worksheets.Item("Silver Busbar front-side")
$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "7,816"

$ws3 = $wb.Worksheets.Item("Silver finger front-side")
$ws3.Range("B10").NumberFormat = "@"
$ws3.Range("B10").Value = "7,866"

$ws4 = $wb.Worksheets.Item("USD_CNY")
$ws4.Range("B10").NumberFormat = "@"
$ws4.Range("B10").Value = "7.2787"
